$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.73%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.78%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.155"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.67%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05623"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.62%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.475"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.26%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8168"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.35%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8315"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.90%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1327"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.86%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06894"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.76%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.02884"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.32%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09382"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.22%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.001509"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.18%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").Value = "'0.04225"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-9.85%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005986"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.38%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.42%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.606"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.50%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.021"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.03%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.226"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'5.10%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3113"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.38%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C21").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D21").Value = "'0.03091"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-3.83%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1291"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.16%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D23").Value = "'3.737"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.04%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.13%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001226"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.93%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004488"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-2.61%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009793"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'2.03%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001383"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-0.50%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03647"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.10%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1367"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.41%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002598"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.30%"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003425"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-43.94%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008175"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.69%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'-0.01%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.05%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-36.77%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002632"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'28.50%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.05%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E50").Style = "Normal"
